$wb = $excel.ActiveWorkbook

# --- Rename the first sheet (was "INTER_SWITCH_LINKS") to "SWITCH_TO_SWITCH" ---
$wsSwitch = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$wsSwitch.Name = "SWITCH_TO_SWITCH"

# --- The author's session ended on SWITCH_TO_SWITCH instead of COMPUTE_NODES ---
# Make it the active/selected sheet (this also clears COMPUTE_NODES's
# tabSelected flag, since only one sheet can be tab-selected at a time).
$wsSwitch.Activate()

# Keep the sheet scrolled so row 4 / column B is the top-left visible cell,
# matching the saved view state, then move the active cell/selection from
# T24 to E29 as recorded in the workbook.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$wsSwitch.Range("E29").Select()
